$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# --- Insert new header columns AA:AE (push old "Status as of ..." header to AF) ---
# First stash the current AA1 value/content before we add new columns in its place.
$oldAA1 = $ws.Range("AA1").Value2

# Copy the formatting of an existing bold/bordered header cell (Z1) into the
# new header cells AA1:AE1 so they pick up the same header style used by the
# rest of row 1.
$ws.Range("Z1").Copy()
$ws.Range("AA1:AE1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"

# Move the old AA1 header ("Status as of July 4, 2025") to the new AF1 cell,
# keeping its original (unstyled) formatting.
$ws.Range("AF1").Value = $oldAA1

# --- Clear placeholder "-" values from columns I and L for rows 2-17 ---
$ws.Range("I2:I17").ClearContents()
$ws.Range("L2:L17").ClearContents()

# --- Add contractor name into newly available AA9 cell ---
$ws.Range("AA9").Value = "URBANCON Builders & Supply"

# --- Update the data validation (dropdown) range from AA2:AA85 to AF2:AF85 ---
$ws.Range("AA2:AA85").Validation.Delete()
$ws.Range("AF2:AF85").Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
$ws.Range("AF2:AF85").Validation.IgnoreBlank = $true
$ws.Range("AF2:AF85").Validation.InCellDropdown = $true
$ws.Range("AF2:AF85").Validation.ShowInput = $false
$ws.Range("AF2:AF85").Validation.ShowError = $false

Write-Host "Edit complete"
